# Applies the changes described by the diff:
#  - Row 19 and Row 20 swap which worker (ID + Name) they show:
#      Row19 (was 73000489 / MIGUEL ENRIQUE ESCOBAR LANDERO) -> becomes 78754019 / ROBERTO CARLOS PONCE ARRAZOLA
#      Row20 (was 78754019 / ROBERTO CARLOS PONCE ARRAZOLA)  -> becomes 73000489 / MIGUEL ENRIQUE ESCOBAR LANDERO
#  - "Valor Mora" (column G) amounts updated on rows 17, 18 and 19:
#      G17: 1181084 -> 1300000
#      G18: 1181084 -> 1508000
#      G19: 1181084 -> 877803  (this is now Roberto's row after the swap above)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the worker identity shown on rows 19 and 20 (column C = ID doc, column D = name)
$ws.Range("C19").Value = "78754019"
$ws.Range("D19").Value = "ROBERTO CARLOS PONCE ARRAZOLA"
$ws.Range("C20").Value = "73000489"
$ws.Range("D20").Value = "MIGUEL ENRIQUE ESCOBAR LANDERO"

# Update "Valor Mora" values in column G
$ws.Range("G17").Value = 1300000
$ws.Range("G18").Value = 1508000
$ws.Range("G19").Value = 877803
